$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header casing: "Jansen_Omschrijving" -> "Jansen_omschrijving"
$ws.Range("C1").Value = "Jansen_omschrijving"

# Replace the shared "=B#" formulas in column C (rows 43-57, 60) with the
# short Dutch category labels used elsewhere in the column, as plain values
# (no longer formulas).
$ws.Range("C43").Value = "Hypofyse- en hypothalamus en verwante verbindingen"
$ws.Range("C44").Value = "Corticosteroiden systemisch"
$ws.Range("C45").Value = "Schildklierhormonen"
$ws.Range("C46").Value = "Pancreashormonen"
$ws.Range("C47").Value = "Calciumregulerende middelen"
$ws.Range("C49").Value = "Antimycotica"
$ws.Range("C50").Value = "Antimycobacteriele middelen"
$ws.Range("C51").Value = "Antivirale middelen"
$ws.Range("C52").Value = "Sera en immunoglobulinen"
$ws.Range("C53").Value = "Vaccins"
$ws.Range("C54").Value = "Oncolytica"
$ws.Range("C55").Value = "Hormonen"
$ws.Range("C56").Value = "Immunostimulantia"
$ws.Range("C57").Value = "Immunosuppressiva"
$ws.Range("C60").Value = "Spierrelaxantia"

# Update the active selection to match the edited workbook
$null = $ws.Range("C2").Select()
